# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for rows 2-30 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 5
    4  = 4
    5  = 1
    6  = 1
    7  = 4
    8  = 2
    9  = 5
    10 = 6
    11 = 6
    12 = 2
    13 = 3
    15 = 0
    16 = 1
    17 = 1
    18 = 3
    19 = 5
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 3
    25 = 1
    27 = 4
    28 = 5
    29 = 2
    30 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
